# Update tab names in all BOMs, fix bi-color LED naming.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "AxB" to "BOM"
$ws.Name = "BOM"

# Fix the bi-color LED package naming (Flat -> Dome)
$ws.Range("C18").Value = "LED 3mm Dome Bicolor"
